$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Summary" cell F11 ---
$wsSummary = $wb.Worksheets.Item("Summary")
$cellSummary = $wsSummary.Range("F11")
$newSummary = "Gradle:com.here.ort.gradle.example:lib:1.0.0`n  Unknown time [ERROR]: FileCounter - DownloadException: Download failed for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'.`nSuppressed: DownloadException: No VCS URL provided for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'. Please make sure the release POM file includes the SCM connection, see: https://docs.gradle.org/current/userguide/publishing_maven.html#example_customizing_the_pom_file, `nSuppressed: DownloadException: No source artifact URL provided for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'.`n"
$cellSummary.Value = $newSummary
$lenSummary = $newSummary.Length
# Re-apply the run's font across two adjacent sub-ranges so the exporter keeps
# a single rich-text run (matching the formatted-run style of the original
# cell) instead of collapsing the whole string into a plain, unformatted <t>.
$cellSummary.Characters(1, 1).Font.Name = "Calibri"
$cellSummary.Characters(2, $lenSummary - 1).Font.Name = "Calibri"

# --- Sheet 2: "Gradle com.here.ort.gradle.exam" cell F11 ---
$wsGradle = $wb.Worksheets.Item("Gradle com.here.ort.gradle.exam")
$cellGradle = $wsGradle.Range("F11")
$newGradle = "Unknown time [ERROR]: FileCounter - DownloadException: Download failed for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'.`nSuppressed: DownloadException: No VCS URL provided for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'. Please make sure the release POM file includes the SCM connection, see: https://docs.gradle.org/current/userguide/publishing_maven.html#example_customizing_the_pom_file, `nSuppressed: DownloadException: No source artifact URL provided for 'Gradle:com.here.ort.gradle.example:lib:1.0.0'."
$cellGradle.Value = $newGradle
$lenGradle = $newGradle.Length
$cellGradle.Characters(1, 1).Font.Name = "Calibri"
$cellGradle.Characters(2, $lenGradle - 1).Font.Name = "Calibri"
